$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF ("Date") rows 2-31 were stamped with the wrong text
# "5-16-2013-14". Because NBA box-score dates were captured off by a
# day, the correct value is "2014-05-16" for every row on this sheet.
# Force the cells to keep storing plain text (not an auto-converted
# date serial) so the corrected value round-trips as literal text.
$rng = $ws.Range("BF2:BF31")
$rng.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)
    if ($cell.Value2 -eq "5-16-2013-14") {
        $cell.Value = "2014-05-16"
    }
}
